# Auto-generated Excel COM-interop script applying numeric updates
# to the Ridill_Profits workbook (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 290.8
$ws.Range("J17").Value = 290.8
$ws.Range("L17").Value = 872.4000000000001
$ws.Range("N17").Value = -1208.4

$ws.Range("H31").Value = 1019.8
$ws.Range("I31").Value = 499.66666
$ws.Range("J31").Value = 1800
$ws.Range("K31").Value = 1498.99998
$ws.Range("L31").Value = 5400
$ws.Range("M31").Value = -1268.99998
$ws.Range("N31").Value = -5860

$ws.Range("H118").Value = 4611.7646
$ws.Range("I118").Value = 500
$ws.Range("J118").Value = 5320.6895
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 15962.0685
$ws.Range("M118").Value = 157
$ws.Range("N118").Value = -19276.0685

$ws.Range("H138").Value = 5148.7334
$ws.Range("I138").Value = 7368.5713
$ws.Range("J138").Value = 3206.375
$ws.Range("K138").Value = 22105.7139
$ws.Range("L138").Value = 9619.125
$ws.Range("M138").Value = -16965.7139
$ws.Range("N138").Value = -19899.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8482107
$ws.Range("I32").Value = 7690.3335
$ws.Range("K32").Value = 7690.3335
$ws.Range("M32").Value = -7403.3335

$ws.Range("H74").Value = 53989172
$ws.Range("I74").Value = 115742600
$ws.Range("J74").Value = 14290542
$ws.Range("K74").Value = 115742600
$ws.Range("L74").Value = 14290542
$ws.Range("M74").Value = -115741726
$ws.Range("N74").Value = -14292290

$ws.Range("H77").Value = 53989172
$ws.Range("I77").Value = 115742600
$ws.Range("J77").Value = 14290542
$ws.Range("K77").Value = 578713000
$ws.Range("L77").Value = 71452710
$ws.Range("M77").Value = -578708632
$ws.Range("N77").Value = -71461446

$ws.Range("H88").Value = 4376.923
$ws.Range("I88").Value = 1983.3334
$ws.Range("J88").Value = 6428.5713
$ws.Range("K88").Value = 1983.3334
$ws.Range("L88").Value = 6428.5713
$ws.Range("M88").Value = -1577.3334
$ws.Range("N88").Value = -7240.5713

$ws.Range("H91").Value = 4376.923
$ws.Range("I91").Value = 1983.3334
$ws.Range("J91").Value = 6428.5713
$ws.Range("K91").Value = 1983.3334
$ws.Range("L91").Value = 6428.5713
$ws.Range("M91").Value = -579.3334
$ws.Range("N91").Value = -9236.5713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1986.06
$ws.Range("I86").Value = 1986.06
$ws.Range("K86").Value = 1986.06
$ws.Range("M86").Value = -863.0599999999999

$ws.Range("H89").Value = 1986.06
$ws.Range("I89").Value = 1986.06
$ws.Range("K89").Value = 9930.299999999999
$ws.Range("M89").Value = -4314.299999999999

$ws.Range("H107").Value = 2503502.8
$ws.Range("I107").Value = 2503502.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2503502.8
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2501582.8
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1633.3334
$ws.Range("I16").Value = 1950
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1950
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1663
$ws.Range("N16").Value = -1574

$ws.Range("H22").Value = 367.15384
$ws.Range("I22").Value = 308.22223
$ws.Range("J22").Value = 499.75
$ws.Range("K22").Value = 308.22223
$ws.Range("L22").Value = 499.75
$ws.Range("M22").Value = 41.77776999999998
$ws.Range("N22").Value = -1199.75

$ws.Range("H41").Value = 5910
$ws.Range("I41").Value = 5910
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5910
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -5482
$ws.Range("N41").ClearContents()

$ws.Range("H51").Value = 19749.5
$ws.Range("J51").Value = 19749.5
$ws.Range("L51").Value = 19749.5
$ws.Range("N51").Value = -21221.5

$ws.Range("H59").Value = 26466
$ws.Range("J59").Value = 26466
$ws.Range("L59").Value = 26466
$ws.Range("N59").Value = -28756

$ws.Range("H60").Value = 5995.3335
$ws.Range("I60").Value = 5995.3335
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 5995.3335
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -5484.3335
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 19749.5
$ws.Range("J61").Value = 19749.5
$ws.Range("L61").Value = 19749.5
$ws.Range("N61").Value = -20445.5

$ws.Range("H68").Value = 19960
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 24200
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 24200
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -25698

$ws.Range("H71").Value = 19960
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 24200
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 72600
$ws.Range("M71").Value = -5256
$ws.Range("N71").Value = -80088

$ws.Range("H74").Value = 35500
$ws.Range("J74").Value = 35500
$ws.Range("L74").Value = 35500
$ws.Range("N74").Value = -37248

$ws.Range("H77").Value = 35500
$ws.Range("J77").Value = 35500
$ws.Range("L77").Value = 106500
$ws.Range("N77").Value = -115236

$ws.Range("H107").Value = 414.45947
$ws.Range("I107").Value = 261
$ws.Range("J107").Value = 734.1667
$ws.Range("K107").Value = 261
$ws.Range("L107").Value = 734.1667
$ws.Range("M107").Value = 1659
$ws.Range("N107").Value = -4574.1667

$ws.Range("H113").Value = 1633.3334
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -5340

$ws.Range("H132").Value = 1472375.1
$ws.Range("I132").Value = 2000957.9
$ws.Range("J132").Value = 4089.5557
$ws.Range("K132").Value = 6002873.699999999
$ws.Range("L132").Value = 12268.6671
$ws.Range("M132").Value = -6000343.699999999
$ws.Range("N132").Value = -17328.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3521.2144
$ws.Range("J68").Value = 5435.88
$ws.Range("L68").Value = 16307.64
$ws.Range("N68").Value = -17929.64

$ws.Range("H71").Value = 3521.2144
$ws.Range("J71").Value = 5435.88
$ws.Range("L71").Value = 48922.92
$ws.Range("N71").Value = -57034.92

$ws.Range("H132").Value = 1746.7646
$ws.Range("I132").Value = 1772.5
$ws.Range("J132").Value = 1738.8462
$ws.Range("K132").Value = 15952.5
$ws.Range("L132").Value = 15649.6158
$ws.Range("M132").Value = -13422.5
$ws.Range("N132").Value = -20709.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 29111.111
$ws.Range("I113").Value = 1460
$ws.Range("J113").Value = 63675
$ws.Range("K113").Value = 1460
$ws.Range("L113").Value = 63675
$ws.Range("M113").Value = 710
$ws.Range("N113").Value = -68015

$ws.Range("H132").Value = 22085690
$ws.Range("I132").Value = 20910274
$ws.Range("J132").Value = 23378650
$ws.Range("K132").Value = 62730822
$ws.Range("L132").Value = 70135950
$ws.Range("M132").Value = -62728292
$ws.Range("N132").Value = -70141010

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1497.5
$ws.Range("I100").Value = 1137.5
$ws.Range("J100").Value = 2217.5
$ws.Range("K100").Value = 1137.5
$ws.Range("L100").Value = 2217.5
$ws.Range("M100").Value = -596.5
$ws.Range("N100").Value = -3299.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 353764.78
$ws.Range("I132").Value = 488981.12
$ws.Range("J132").Value = 2202.3
$ws.Range("K132").Value = 2202.3
$ws.Range("L132").Value = 6606.900000000001
$ws.Range("M132").Value = -1464413.36
$ws.Range("N132").Value = -11666.9
